# Auto-generated COM-interop script applying scheduled-runner market-data refresh
# to the Leve profit sheets (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit*).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 969.67645
$ws.Range("J19").Value = 1068.4736
$ws.Range("L19").Value = 1068.4736
$ws.Range("N19").Value = -1418.4736
$ws.Range("H28").Value = 2149
$ws.Range("I28").Value = 2118.88
$ws.Range("K28").Value = 2118.88
$ws.Range("M28").Value = -1633.88
$ws.Range("H76").Value = 4333
$ws.Range("J76").Value = 4999.5
$ws.Range("L76").Value = 4999.5
$ws.Range("N76").Value = -5629.5
$ws.Range("H79").Value = 4333
$ws.Range("J79").Value = 4999.5
$ws.Range("L79").Value = 4999.5
$ws.Range("N79").Value = -7183.5
$ws.Range("H93").Value = 47494.5
$ws.Range("J93").Value = 47494.5
$ws.Range("L93").Value = 47494.5
$ws.Range("N93").Value = -52486.5
$ws.Range("H132").Value = 1991.2131
$ws.Range("I132").Value = 1990.4642
$ws.Range("J132").Value = 1999.6
$ws.Range("K132").Value = 5971.392599999999
$ws.Range("L132").Value = 5998.799999999999
$ws.Range("M132").Value = -3441.392599999999
$ws.Range("N132").Value = -11058.8
$ws.Range("H135").Value = 851.9535
$ws.Range("I135").Value = 730.8421
$ws.Range("J135").Value = 1772.4
$ws.Range("K135").Value = 6577.5789
$ws.Range("L135").Value = 15951.6
$ws.Range("M135").Value = -4042.5789
$ws.Range("N135").Value = -21021.6
$ws.Range("H138").Value = 11908512
$ws.Range("I138").Value = 1621.5
$ws.Range("J138").Value = 13892993
$ws.Range("K138").Value = 4864.5
$ws.Range("L138").Value = 41678979
$ws.Range("M138").Value = 275.5
$ws.Range("N138").Value = -41689259

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 25000
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H54").Value = 50000
$ws.Range("I54").Value = 50000
$ws.Range("K54").Value = 50000
$ws.Range("M54").Value = -49231
$ws.Range("H109").Value = 42099.332
$ws.Range("J109").Value = 42099.332
$ws.Range("L109").Value = 42099.332
$ws.Range("N109").Value = -44873.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5055.36
$ws.Range("I86").Value = 5005.5
$ws.Range("K86").Value = 5005.5
$ws.Range("M86").Value = -3882.5
$ws.Range("H89").Value = 5055.36
$ws.Range("I89").Value = 5005.5
$ws.Range("K89").Value = 25027.5
$ws.Range("M89").Value = -19411.5
$ws.Range("H99").Value = 63581
$ws.Range("I99").Value = 45105.348
$ws.Range("K99").Value = 45105.348
$ws.Range("M99").Value = -43607.348
$ws.Range("H134").Value = 1673.7805
$ws.Range("I134").Value = 1456.1842
$ws.Range("K134").Value = 4368.5526
$ws.Range("M134").Value = -1833.5526

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2870.4443
$ws.Range("J94").Value = 3400.2
$ws.Range("L94").Value = 3400.2
$ws.Range("N94").Value = -4302.2
$ws.Range("H107").Value = 600.4375
$ws.Range("I107").Value = 577.46155
$ws.Range("K107").Value = 577.46155
$ws.Range("M107").Value = 1342.53845
$ws.Range("H132").Value = 6182.6
$ws.Range("I132").Value = 4999.6665
$ws.Range("K132").Value = 14998.9995
$ws.Range("M132").Value = -12468.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 847.65
$ws.Range("I5").Value = 816.35297
$ws.Range("K5").Value = 2449.05891
$ws.Range("M5").Value = -2337.05891
$ws.Range("H8").Value = 999
$ws.Range("I8").Value = 999
$ws.Range("K8").Value = 2997
$ws.Range("M8").Value = -2858
$ws.Range("H29").Value = 382.2857
$ws.Range("I29").Value = 170.5
$ws.Range("J29").Value = 467
$ws.Range("K29").Value = 511.5
$ws.Range("L29").Value = 1401
$ws.Range("M29").Value = -234.5
$ws.Range("N29").Value = -1955
$ws.Range("H68").Value = 8334079.5
$ws.Range("J68").Value = 12500845
$ws.Range("L68").Value = 37502535
$ws.Range("N68").Value = -37504157
$ws.Range("H71").Value = 8334079.5
$ws.Range("J71").Value = 12500845
$ws.Range("L71").Value = 112507605
$ws.Range("N71").Value = -112515717
$ws.Range("H122").Value = 993.087
$ws.Range("I122").Value = 570.38464
$ws.Range("J122").Value = 1542.6
$ws.Range("K122").Value = 5133.46176
$ws.Range("L122").Value = 13883.4
$ws.Range("M122").Value = -2683.46176
$ws.Range("N122").Value = -18783.4
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("H133").Value = 16553.691
$ws.Range("J133").Value = 11714.286
$ws.Range("L133").Value = 35142.858
$ws.Range("N133").Value = -45262.858
$ws.Range("H135").Value = 847.65
$ws.Range("I135").Value = 816.35297
$ws.Range("K135").Value = 7347.17673
$ws.Range("M135").Value = -4812.17673
$ws.Range("H136").Value = 675403.9399999999
$ws.Range("I136").Value = 1115784.4
$ws.Range("K136").Value = 3347353.2
$ws.Range("M136").Value = -3342253.2
$ws.Range("H137").Value = 2887.3333
$ws.Range("I137").Value = 1563.8572
$ws.Range("K137").Value = 4691.571599999999
$ws.Range("M137").Value = 408.4284000000007
$ws.Range("H138").Value = 29421274
$ws.Range("I138").Value = 35722620
$ws.Range("K138").Value = 107167860
$ws.Range("M138").Value = -107162720
$ws.Range("H139").Value = 3583.1667
$ws.Range("I139").Value = 3924.75
$ws.Range("J139").Value = 2900
$ws.Range("K139").Value = 11774.25
$ws.Range("L139").Value = 8700
$ws.Range("M139").Value = -6634.25
$ws.Range("N139").Value = -18980
$ws.Range("H140").Value = 1716
$ws.Range("I140").Value = 1716
$ws.Range("K140").Value = 5148
$ws.Range("M140").Value = 32

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3101.5
$ws.Range("I80").Value = 3120.8
$ws.Range("K80").Value = 3120.8
$ws.Range("M80").Value = -2122.8
$ws.Range("H83").Value = 3101.5
$ws.Range("I83").Value = 3120.8
$ws.Range("K83").Value = 15604
$ws.Range("M83").Value = -10612
$ws.Range("H126").Value = 11838.292
$ws.Range("I126").Value = 13467.7
$ws.Range("J126").Value = 3691.25
$ws.Range("K126").Value = 40403.10000000001
$ws.Range("L126").Value = 11073.75
$ws.Range("M126").Value = -37933.10000000001
$ws.Range("N126").Value = -16013.75
$ws.Range("H132").Value = 3827.5
$ws.Range("I132").Value = 3419.4443
$ws.Range("K132").Value = 10258.3329
$ws.Range("M132").Value = -7728.332900000001
$ws.Range("H136").Value = 36536.957
$ws.Range("J136").Value = 36536.957
$ws.Range("L136").Value = 109610.871
$ws.Range("N136").Value = -114710.871

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2231.8
$ws.Range("I16").Value = 2554.5789
$ws.Range("J16").Value = 1209.6666
$ws.Range("K16").Value = 2554.5789
$ws.Range("L16").Value = 1209.6666
$ws.Range("M16").Value = -2384.5789
$ws.Range("N16").Value = -1549.6666
$ws.Range("H40").Value = 4258.6875
$ws.Range("I40").Value = 3289.9092
$ws.Range("K40").Value = 3289.9092
$ws.Range("M40").Value = -3153.9092
$ws.Range("H76").Value = 25000
$ws.Range("J76").Value = 25000
$ws.Range("L76").Value = 25000
$ws.Range("N76").Value = -25676
$ws.Range("H79").Value = 25000
$ws.Range("J79").Value = 25000
$ws.Range("L79").Value = 25000
$ws.Range("N79").Value = -27340
$ws.Range("H132").Value = 24484.541
$ws.Range("I132").Value = 34607.062
$ws.Range("J132").Value = 4239.5
$ws.Range("K132").Value = 103821.186
$ws.Range("L132").Value = 12718.5
$ws.Range("M132").Value = -101291.186
$ws.Range("N132").Value = -17778.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5958.433
$ws.Range("I81").Value = 8969.076999999999
$ws.Range("K81").Value = 17938.154
$ws.Range("M81").Value = -16877.154
$ws.Range("H84").Value = 5958.433
$ws.Range("I84").Value = 8969.076999999999
$ws.Range("K84").Value = 89690.76999999999
$ws.Range("M84").Value = -84386.76999999999
$ws.Range("H107").Value = 33388.516
$ws.Range("I107").Value = 1220.125
$ws.Range("J107").Value = 143680.14
$ws.Range("K107").Value = 3660.375
$ws.Range("L107").Value = 431040.42
$ws.Range("M107").Value = -1740.375
$ws.Range("N107").Value = -434880.42
$ws.Range("H132").Value = 2536.6316
$ws.Range("I132").Value = 2365.7812
$ws.Range("J132").Value = 3447.8333
$ws.Range("K132").Value = 7097.3436
$ws.Range("L132").Value = 10343.4999
$ws.Range("M132").Value = -4567.3436
$ws.Range("N132").Value = -15403.4999
